# Updates the cryptocurrency price/volume table (rows 2-51) on the active
# sheet with refreshed market data, matching the upstream GitHub Actions
# data-refresh commit. Coin/Link (columns B/C) are only rewritten where the
# ranking order changed; Price (D) and Volume(1h) (E) are refreshed for the
# affected rows.
#
# Values are written with a leading "'" so Excel stores them as text
# (preventing numeric-looking prices like "1.00"/"2.30" from being
# normalised to "1"/"2.3", and preserving the multi-dot thousands format
# like "36.647.69"). The style is reset to "Normal" afterwards so the
# quote-prefix formatting flag introduced by the text-forcing trick does
# not leave a stray cell style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.647.69"
$ws.Range("E2").Value = "'  -1.06%  "
$ws.Range("D3").Value = "'2.055.56"
$ws.Range("E3").Value = "'  +0.01%  "
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("D5").Value = "'248.26"
$ws.Range("E5").Value = "'  +0.73%  "
$ws.Range("D6").Value = "'0.664"
$ws.Range("E6").Value = "'  +1.04%  "
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("D8").Value = "'54.76"
$ws.Range("E8").Value = "'  -7.02%  "
$ws.Range("D9").Value = "'60.72"
$ws.Range("E9").Value = "'  +1.95%  "
$ws.Range("D10").Value = "'0.367"
$ws.Range("E10").Value = "'  -2.79%  "
$ws.Range("D11").Value = "'0.0755"
$ws.Range("E11").Value = "'  -2.42%  "
$ws.Range("E12").Value = "'  -3.04%  "
$ws.Range("D13").Value = "'0.973"
$ws.Range("E13").Value = "'  +9.90%  "
$ws.Range("D14").Value = "'14.82"
$ws.Range("E14").Value = "'  -4.22%  "
$ws.Range("D15").Value = "'2.358.68"
$ws.Range("E15").Value = "'  +0.15%  "
$ws.Range("E16").Value = "'  -4.17%  "
$ws.Range("D17").Value = "'2.059.48"
$ws.Range("E17").Value = "'  +1.20%  "
$ws.Range("D18").Value = "'36.537.45"
$ws.Range("E18").Value = "'  -1.24%  "
$ws.Range("D19").Value = "'17.35"
$ws.Range("E19").Value = "'  -4.49%  "
$ws.Range("D20").Value = "'72.18"
$ws.Range("E20").Value = "'  -2.44%  "
$ws.Range("D21").Value = "'0.0₃0862"
$ws.Range("E21").Value = "'  -3.17%  "
$ws.Range("D22").Value = "'238.50"
$ws.Range("E22").Value = "'  +0.14%  "
$ws.Range("E23").Value = "'  -3.30%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "'  +0.02%  "
$ws.Range("E25").Value = "'  -2.51%  "
$ws.Range("D26").Value = "'2.30"
$ws.Range("E26").Value = "'  +6.47%  "
$ws.Range("D27").Value = "'166.36"
$ws.Range("E27").Value = "'  -1.90%  "
$ws.Range("D28").Value = "'9.26"
$ws.Range("E28").Value = "'  -8.06%  "
$ws.Range("D29").Value = "'20.14"
$ws.Range("E29").Value = "'  +0.16%  "
$ws.Range("E30").Value = "'  -1.62%  "
$ws.Range("E31").Value = "'  +9.37%  "
$ws.Range("E32").Value = "'  -6.23%  "
$ws.Range("D33").Value = "'4.51"
$ws.Range("E33").Value = "'  -3.70%  "
$ws.Range("D34").Value = "'0.0596"
$ws.Range("E34").Value = "'  -3.44%  "
$ws.Range("E35").Value = "'  +0.10%  "
$ws.Range("D36").Value = "'0.0865"
$ws.Range("E36").Value = "'  +2.71%  "
$ws.Range("D37").Value = "'2.27"
$ws.Range("E37").Value = "'  -1.90%  "
$ws.Range("D38").Value = "'1.83"
$ws.Range("E38").Value = "'  -0.08%  "
$ws.Range("B39").Value = "'TrustWalletToken"
$ws.Range("C39").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.26"
$ws.Range("E39").Value = "'  -5.27%  "
$ws.Range("B40").Value = "'THORChain"
$ws.Range("C40").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").Value = "'5.04"
$ws.Range("E40").Value = "'  -3.87%  "
$ws.Range("E41").Value = "'  -5.41%  "
$ws.Range("E42").Value = "'  -3.84%  "
$ws.Range("B43").Value = "'ARBITRUM"
$ws.Range("C43").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.11"
$ws.Range("E43").Value = "'  -4.86%  "
$ws.Range("B44").Value = "'Aave"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'94.93"
$ws.Range("E44").Value = "'  -2.87%  "
$ws.Range("D45").Value = "'0.0918"
$ws.Range("E45").Value = "'  -4.60%  "
$ws.Range("D46").Value = "'1.415.06"
$ws.Range("E46").Value = "'  +8.64%  "
$ws.Range("B47").Value = "'InjectiveProtocol"
$ws.Range("C47").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'15.99"
$ws.Range("E47").Value = "'  -5.93%  "
$ws.Range("B48").Value = "'FraxShare"
$ws.Range("C48").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.58"
$ws.Range("E48").Value = "'  +11.14%  "
$ws.Range("D49").Value = "'2.92"
$ws.Range("E49").Value = "'  +1.80%  "
$ws.Range("E50").Value = "'  -3.78%  "
$ws.Range("D51").Value = "'46.05"
$ws.Range("E51").Value = "'  +3.66%  "

# Clear the text quote-prefix styling introduced above so the cells keep
# their original (default/"Normal") style.
$ws.Range("B2:E51").Style = "Normal"
